$d = $word.ActiveDocument
$frag = @'
    <w:p>
      <w:pPr>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t xml:space="preserve">TS Pada Paatam – TS 2.6 </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>Sanskrit</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t xml:space="preserve"> co</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>rrections –</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t xml:space="preserve">Observed </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t xml:space="preserve">till </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:highlight w:val="red"/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>????</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t xml:space="preserve"> (</w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>ignore</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t xml:space="preserve"> those which are already incorporated in your book’s version and date)</w:t>
      </w:r>
    </w:p>
    <w:tbl>
      <w:tblPr>
        <w:tblW w:w="12794" w:type="dxa"/>
        <w:tblInd w:w="-318" w:type="dxa"/>
        <w:tblBorders>
          <w:top w:val="single" w:sz="4" w:space="0" w:color="000000"/>
          <w:left w:val="single" w:sz="4" w:space="0" w:color="000000"/>
          <w:bottom w:val="single" w:sz="4" w:space="0" w:color="000000"/>
          <w:right w:val="single" w:sz="4" w:space="0" w:color="000000"/>
          <w:insideH w:val="single" w:sz="4" w:space="0" w:color="000000"/>
          <w:insideV w:val="single" w:sz="4" w:space="0" w:color="000000"/>
        </w:tblBorders>
        <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
      </w:tblPr>
      <w:tblGrid>
        <w:gridCol w:w="3120"/>
        <w:gridCol w:w="4424"/>
        <w:gridCol w:w="5250"/>
      </w:tblGrid>
      <w:tr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3120" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
              </w:rPr>
              <w:t>Section, Paragraph</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
              </w:rPr>
              <w:t>Reference</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="4424" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
              </w:rPr>
              <w:t>As Printed</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="5250" w:type="dxa"/>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
              <w:ind w:right="-18"/>
              <w:jc w:val="center"/>
              <w:rPr>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="32"/>
                <w:szCs w:val="32"/>
              </w:rPr>
              <w:t>To be read as or corrected as</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
      <w:tr>
        <w:trPr>
          <w:trHeight w:val="569"/>
        </w:trPr>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="3120" w:type="dxa"/>
            <w:tcBorders>
              <w:top w:val="single" w:sz="4" w:space="0" w:color="000000"/>
              <w:left w:val="single" w:sz="4" w:space="0" w:color="000000"/>
              <w:bottom w:val="single" w:sz="4" w:space="0" w:color="000000"/>
              <w:right w:val="single" w:sz="4" w:space="0" w:color="000000"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:widowControl w:val="0"/>
              <w:autoSpaceDE w:val="0"/>
              <w:autoSpaceDN w:val="0"/>
              <w:adjustRightInd w:val="0"/>
              <w:spacing w:before="0" w:line="264" w:lineRule="auto"/>
              <w:ind w:right="-320"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:eastAsia="en-IN" w:bidi="ml-IN"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:eastAsia="en-IN" w:bidi="ml-IN"/>
              </w:rPr>
              <w:t>TS 2.6.</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:eastAsia="en-IN" w:bidi="ml-IN"/>
              </w:rPr>
              <w:t>10</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:eastAsia="en-IN" w:bidi="ml-IN"/>
              </w:rPr>
              <w:t>.</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:eastAsia="en-IN" w:bidi="ml-IN"/>
              </w:rPr>
              <w:t>1</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:eastAsia="en-IN" w:bidi="ml-IN"/>
              </w:rPr>
              <w:t xml:space="preserve"> - Vaakyam</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:widowControl w:val="0"/>
              <w:autoSpaceDE w:val="0"/>
              <w:autoSpaceDN w:val="0"/>
              <w:adjustRightInd w:val="0"/>
              <w:spacing w:before="0" w:line="264" w:lineRule="auto"/>
              <w:ind w:right="-320"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:eastAsia="en-IN" w:bidi="ml-IN"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:eastAsia="en-IN" w:bidi="ml-IN"/>
              </w:rPr>
              <w:t xml:space="preserve">Line No. </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:eastAsia="en-IN" w:bidi="ml-IN"/>
              </w:rPr>
              <w:t>4</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:widowControl w:val="0"/>
              <w:autoSpaceDE w:val="0"/>
              <w:autoSpaceDN w:val="0"/>
              <w:adjustRightInd w:val="0"/>
              <w:spacing w:before="0" w:line="264" w:lineRule="auto"/>
              <w:ind w:right="-320"/>
              <w:jc w:val="both"/>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:eastAsia="en-IN" w:bidi="ml-IN"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:eastAsia="en-IN" w:bidi="ml-IN"/>
              </w:rPr>
              <w:t xml:space="preserve">Panchaati No. </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:eastAsia="en-IN" w:bidi="ml-IN"/>
              </w:rPr>
              <w:t>57</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="4424" w:type="dxa"/>
            <w:tcBorders>
              <w:top w:val="single" w:sz="4" w:space="0" w:color="000000"/>
              <w:left w:val="single" w:sz="4" w:space="0" w:color="000000"/>
              <w:bottom w:val="single" w:sz="4" w:space="0" w:color="000000"/>
              <w:right w:val="single" w:sz="4" w:space="0" w:color="000000"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:widowControl w:val="0"/>
              <w:autoSpaceDE w:val="0"/>
              <w:autoSpaceDN w:val="0"/>
              <w:adjustRightInd w:val="0"/>
              <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
              <w:ind w:right="-23"/>
              <w:rPr>
                <w:rFonts w:ascii="BRH Malayalam Extra" w:hAnsi="BRH Malayalam Extra" w:cs="BRH Malayalam Extra"/>
                <w:sz w:val="40"/>
                <w:szCs w:val="40"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Devanagari Extra" w:hAnsi="BRH Devanagari Extra" w:cs="BRH Devanagari Extra"/>
                <w:sz w:val="40"/>
                <w:szCs w:val="40"/>
              </w:rPr>
              <w:t>rÉ</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Devanagari Extra" w:hAnsi="BRH Devanagari Extra" w:cs="BRH Devanagari Extra"/>
                <w:sz w:val="40"/>
                <w:szCs w:val="40"/>
                <w:highlight w:val="yellow"/>
              </w:rPr>
              <w:t>Sç</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Devanagari Extra" w:hAnsi="BRH Devanagari Extra" w:cs="BRH Devanagari Extra"/>
                <w:sz w:val="40"/>
                <w:szCs w:val="40"/>
              </w:rPr>
              <w:t xml:space="preserve"> oÉëÉÿ¼hÉÉå</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Malayalam Extra" w:hAnsi="BRH Malayalam Extra" w:cs="BRH Devanagari Extra"/>
                <w:sz w:val="32"/>
                <w:szCs w:val="40"/>
              </w:rPr>
              <w:t>–</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Devanagari Extra" w:hAnsi="BRH Devanagari Extra" w:cs="BRH Devanagari Extra"/>
                <w:sz w:val="40"/>
                <w:szCs w:val="40"/>
              </w:rPr>
              <w:t>£üÉå Å´Éþ¬kÉÉlÉÉå</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Malayalam Extra" w:hAnsi="BRH Malayalam Extra" w:cs="BRH Devanagari Extra"/>
                <w:sz w:val="32"/>
                <w:szCs w:val="40"/>
              </w:rPr>
              <w:t>–</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:eastAsia="Times New Roman" w:cs="Arial"/>
                <w:b/>
                <w:bCs/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:lang w:eastAsia="en-IN" w:bidi="ml-IN"/>
              </w:rPr>
              <w:t xml:space="preserve"> rÉeÉþiÉå</w:t>
            </w:r>
          </w:p>
        </w:tc>
        <w:tc>
          <w:tcPr>
            <w:tcW w:w="5250" w:type="dxa"/>
            <w:tcBorders>
              <w:top w:val="single" w:sz="4" w:space="0" w:color="000000"/>
              <w:left w:val="single" w:sz="4" w:space="0" w:color="000000"/>
              <w:bottom w:val="single" w:sz="4" w:space="0" w:color="000000"/>
              <w:right w:val="single" w:sz="4" w:space="0" w:color="000000"/>
            </w:tcBorders>
          </w:tcPr>
          <w:p>
            <w:pPr>
              <w:widowControl w:val="0"/>
              <w:autoSpaceDE w:val="0"/>
              <w:autoSpaceDN w:val="0"/>
              <w:adjustRightInd w:val="0"/>
              <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
              <w:ind w:right="-22"/>
              <w:rPr>
                <w:rFonts w:ascii="BRH Devanagari Extra" w:hAnsi="BRH Devanagari Extra" w:cs="BRH Devanagari Extra"/>
                <w:sz w:val="40"/>
                <w:szCs w:val="40"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Devanagari Extra" w:hAnsi="BRH Devanagari Extra" w:cs="BRH Devanagari Extra"/>
                <w:sz w:val="40"/>
                <w:szCs w:val="40"/>
              </w:rPr>
              <w:t>rÉ</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Devanagari Extra" w:hAnsi="BRH Devanagari Extra" w:cs="BRH Devanagari Extra"/>
                <w:sz w:val="40"/>
                <w:szCs w:val="40"/>
                <w:highlight w:val="green"/>
              </w:rPr>
              <w:t>S</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Devanagari Extra" w:hAnsi="BRH Devanagari Extra" w:cs="BRH Devanagari Extra"/>
                <w:sz w:val="40"/>
                <w:szCs w:val="40"/>
              </w:rPr>
              <w:t>oÉëÉÿ¼hÉÉå</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Malayalam Extra" w:hAnsi="BRH Malayalam Extra" w:cs="BRH Devanagari Extra"/>
                <w:sz w:val="32"/>
                <w:szCs w:val="40"/>
              </w:rPr>
              <w:t>–</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Devanagari Extra" w:hAnsi="BRH Devanagari Extra" w:cs="BRH Devanagari Extra"/>
                <w:sz w:val="40"/>
                <w:szCs w:val="40"/>
              </w:rPr>
              <w:t>£üÉå Å´Éþ¬kÉÉlÉÉå</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Malayalam Extra" w:hAnsi="BRH Malayalam Extra" w:cs="BRH Devanagari Extra"/>
                <w:sz w:val="32"/>
                <w:szCs w:val="40"/>
              </w:rPr>
              <w:t>–</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Devanagari Extra" w:hAnsi="BRH Devanagari Extra" w:cs="BRH Devanagari Extra"/>
                <w:sz w:val="40"/>
                <w:szCs w:val="40"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:widowControl w:val="0"/>
              <w:autoSpaceDE w:val="0"/>
              <w:autoSpaceDN w:val="0"/>
              <w:adjustRightInd w:val="0"/>
              <w:spacing w:before="0" w:line="240" w:lineRule="auto"/>
              <w:ind w:right="-22"/>
              <w:rPr>
                <w:rFonts w:ascii="BRH Malayalam Extra" w:hAnsi="BRH Malayalam Extra" w:cs="BRH Malayalam Extra"/>
                <w:sz w:val="40"/>
                <w:szCs w:val="40"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="BRH Devanagari Extra" w:hAnsi="BRH Devanagari Extra" w:cs="BRH Devanagari Extra"/>
                <w:sz w:val="40"/>
                <w:szCs w:val="40"/>
              </w:rPr>
              <w:t>rÉeÉþiÉå</w:t>
            </w:r>
          </w:p>
        </w:tc>
      </w:tr>
    </w:tbl>
    <w:p>
      <w:pPr>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
          <w:sz w:val="32"/>
          <w:szCs w:val="32"/>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
    </w:p>
'@
$r = $d.Range(0, 0)
$r.InsertXML($frag)
